# "Modifications before Humanis Training"
#
# 1) Slide 5 - "Text Box 2": move the box right and re-capitalize
#    "eXo specific configuration" -> "eXo Specific Configuration"
# 2) Slide 6 - title "Rectangle 1": drop the stray trailing <a:endParaRPr>
#    (orange FF9900) left after the text, keep the text/run-split/colour.
# 3) Slide 7 - title "Rectangle 1": same cleanup as slide 6.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 5
# ---------------------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$sh5 = $s5.Shapes.Item(2)

# Shape.Left/.Top are in points; the target OOXML offset is 1440606 EMU
# (1 pt = 12700 EMU). 1440606/12700 rounds to a hair under the boundary in
# single precision, so nudge just enough to land back on 1440606 EMU.
$sh5.Left = 113.43355

$tr5 = $sh5.TextFrame.TextRange
# "eXo specific configuration" -> "eXo Specific Configuration"
$tr5.Characters(1, 4).Text = "eXo "
$tr5.Characters(5, 1).Text = "S"
$tr5.Characters(14, 1).Text = "C"

# ---------------------------------------------------------------------------
# Slide 6
# ---------------------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$sh6 = $s6.Shapes.Item(1)
$tr6 = $sh6.TextFrame.TextRange
$tr6.Delete()
$tr6.InsertAfter("eXo")
$tr6.InsertAfter(" ")
$tr6.InsertAfter("specific")
$tr6.InsertAfter(" configuration")

# ---------------------------------------------------------------------------
# Slide 7
# ---------------------------------------------------------------------------
$s7 = $p.Slides.Item(7)
$sh7 = $s7.Shapes.Item(1)
$tr7 = $sh7.TextFrame.TextRange
$tr7.Delete()
$tr7.InsertAfter("eXo")
$tr7.InsertAfter(" ")
$tr7.InsertAfter("specific")
$tr7.InsertAfter(" configuration")
